$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

Set-TextCell $ws.Range("D2") "69.426.86"
Set-TextCell $ws.Range("E2") "  -1.08%  "

Set-TextCell $ws.Range("D3") "2.523.34"
Set-TextCell $ws.Range("E3") "  +0.08%  "

Set-TextCell $ws.Range("E4") "  +0.01%  "

Set-TextCell $ws.Range("D5") "573.57"
Set-TextCell $ws.Range("E5") "  -0.24%  "

Set-TextCell $ws.Range("D6") "166.96"
Set-TextCell $ws.Range("E6") "  -1.69%  "

Set-TextCell $ws.Range("E7") "  -0.06%  "

Set-TextCell $ws.Range("E8") "  +1.48%  "

Set-TextCell $ws.Range("D9") "2.522.54"
Set-TextCell $ws.Range("E9") "  +0.06%  "

Set-TextCell $ws.Range("D10") "0.162"
Set-TextCell $ws.Range("E10") "  +0.89%  "

Set-TextCell $ws.Range("D12") "0.357"
Set-TextCell $ws.Range("E12") "  +4.00%  "

Set-TextCell $ws.Range("D13") "4.91"
Set-TextCell $ws.Range("E13") "  +2.25%  "

Set-TextCell $ws.Range("D14") "2.987.29"
Set-TextCell $ws.Range("E14") "  -0.10%  "

Set-TextCell $ws.Range("D15") "69.297.41"
Set-TextCell $ws.Range("E15") "  -1.10%  "

Set-TextCell $ws.Range("E16") "  -1.77%  "

Set-TextCell $ws.Range("D17") "24.90"
Set-TextCell $ws.Range("E17") "  +0.00%  "

Set-TextCell $ws.Range("D18") "2.520.77"
Set-TextCell $ws.Range("E18") "  -0.43%  "

Set-TextCell $ws.Range("D19") "11.36"
Set-TextCell $ws.Range("E19") "  -0.46%  "

Set-TextCell $ws.Range("D20") "7.63"
Set-TextCell $ws.Range("E20") "  +1.07%  "

Set-TextCell $ws.Range("D21") "349.26"
Set-TextCell $ws.Range("E21") "  -1.28%  "

Set-TextCell $ws.Range("D22") "3.92"
Set-TextCell $ws.Range("E22") "  -0.20%  "

Set-TextCell $ws.Range("E23") "  +1.68%  "

Set-TextCell $ws.Range("E24") "  -0.03%  "

Set-TextCell $ws.Range("D25") "70.40"
Set-TextCell $ws.Range("E25") "  +2.15%  "

Set-TextCell $ws.Range("D26") "3.97"
Set-TextCell $ws.Range("E26") "  -2.43%  "

Set-TextCell $ws.Range("E27") "  -3.39%  "

Set-TextCell $ws.Range("E28") "  -0.35%  "

Set-TextCell $ws.Range("D29") "0.993"
Set-TextCell $ws.Range("E29") "  -0.66%  "

Set-TextCell $ws.Range("E30") "  -1.39%  "

Set-TextCell $ws.Range("D31") "7.83"
Set-TextCell $ws.Range("E31") "  -0.10%  "

Set-TextCell $ws.Range("D32") "463.72"
Set-TextCell $ws.Range("E32") "  -3.59%  "

Set-TextCell $ws.Range("E34") "  -1.25%  "

Set-TextCell $ws.Range("E35") "  -0.01%  "

Set-TextCell $ws.Range("D37") "157.41"
Set-TextCell $ws.Range("E37") "  +0.17%  "

Set-TextCell $ws.Range("D39") "18.58"
Set-TextCell $ws.Range("E39") "  +0.14%  "

Set-TextCell $ws.Range("E40") "  -0.03%  "

Set-TextCell $ws.Range("B41") "PolygonEcosystemToken"
Set-TextCell $ws.Range("C41") "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextCell $ws.Range("D41") "0.320"
Set-TextCell $ws.Range("E41") "  +0.11%  "

Set-TextCell $ws.Range("B42") "RenderToken"
Set-TextCell $ws.Range("C42") "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextCell $ws.Range("D42") "4.72"
Set-TextCell $ws.Range("E42") "  +0.24%  "

Set-TextCell $ws.Range("E43") "  -2.29%  "

Set-TextCell $ws.Range("D44") "38.13"
Set-TextCell $ws.Range("E44") "  -0.40%  "

Set-TextCell $ws.Range("E45") "  -4.72%  "

Set-TextCell $ws.Range("E46") "  -13.21%  "

Set-TextCell $ws.Range("D47") "141.93"
Set-TextCell $ws.Range("E47") "  -0.16%  "

Set-TextCell $ws.Range("E48") "  +0.06%  "

Set-TextCell $ws.Range("E50") "  -0.14%  "

Set-TextCell $ws.Range("D51") "0.581"
Set-TextCell $ws.Range("E51") "  -2.90%  "
